$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows before row 1087; old rows 1087-1166 shift down to 1089-1168.
$ws.Range("A1087:A1088").EntireRow.Insert()

# New row 1087: Navel / Calibre 88 / EE.UU. import data
$ws.Range("A1087").Value = 10
$ws.Range("B1087").Value = "Vega Modelo de Temuco"
$ws.Range("C1087").Value = "La Araucanía"
$ws.Range("D1087").Value = 44931
$ws.Range("E1087").Value = 9
$ws.Range("F1087").Value = "Fruta"
$ws.Range("G1087").Value = 100102
$ws.Range("H1087").Value = "Cítricos"
$ws.Range("I1087").Value = 100102005
$ws.Range("J1087").Value = "Naranja"
$ws.Range("K1087").Value = "Navel"
$ws.Range("L1087").Value = "Calibre 88"
$ws.Range("M1087").Value = 600
$ws.Range("N1087").Value = 30000
$ws.Range("O1087").Value = 32000
$ws.Range("P1087").Value = 31067
$ws.Range("Q1087").Value = "$/caja 18 kilos importada"
$ws.Range("R1087").Value = "EE.UU."
$ws.Range("S1087").Value = 1726
$ws.Range("T1087").Value = 18

# New row 1088: Navel Late / Primera
$ws.Range("A1088").Value = 10
$ws.Range("B1088").Value = "Vega Modelo de Temuco"
$ws.Range("C1088").Value = "La Araucanía"
$ws.Range("D1088").Value = 44931
$ws.Range("E1088").Value = 9
$ws.Range("F1088").Value = "Fruta"
$ws.Range("G1088").Value = 100102
$ws.Range("H1088").Value = "Cítricos"
$ws.Range("I1088").Value = 100102005
$ws.Range("J1088").Value = "Naranja"
$ws.Range("K1088").Value = "Navel Late"
$ws.Range("L1088").Value = "Primera"
$ws.Range("M1088").Value = 155
$ws.Range("N1088").Value = 12000
$ws.Range("O1088").Value = 12000
$ws.Range("P1088").Value = 12000
$ws.Range("Q1088").Value = "$/bandeja 15 kilos granel"
$ws.Range("R1088").Value = "Región de O'Higgins"
$ws.Range("S1088").Value = 800
$ws.Range("T1088").Value = 15
